# Automatische test-sync: 2025-08-03 18:12:50
#
# Appends the newest test-mail log entry (row 29) to the "Logs" sheet,
# extends the conditional formatting ranges that tracked the previous
# last row (28) so they cover the new last row (29), and refreshes the
# "Dashboard" summary sheet to reflect the updated category counts
# (Planning / Afspraak now leads Intern verzoek / Actie voor medewerker).

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append the new row ----------------------------------
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A29").Value = "Kun jij dit even regelen?"
$ws.Range("B29").Value = "mailmind.test@zohomail.eu"
$ws.Range("C29").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("D29").Value = "Planning / Afspraak"
$ws.Range("E29").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Range("F29").Value = "2025-08-03 18:12:40"
$ws.Range("G29").Value = "Ja"
$ws.Range("H29").Value = "Ja"
$ws.Range("I29").Value = "Nee"
$ws.Range("J29").Value = "Nee"

# ---- Logs sheet: extend conditional formatting to row 29 -------------
$ws.Range("D2:D28").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D29"))
$ws.Range("G2:G28").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G29"))
$ws.Range("H2:H28").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H29"))
$ws.Range("I2:I28").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I29"))
$ws.Range("J2:J28").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J29"))

# ---- Dashboard sheet: refresh the category counts ---------------------
# "Planning / Afspraak" now has 7 hits (was 6) and overtakes
# "Intern verzoek / Actie voor medewerker" (still 6) in the ranking.
$ws2 = $wb.Worksheets.Item("Dashboard")

$ws2.Range("A3").Value = "Planning / Afspraak"
$ws2.Range("B3").Value = 7
$ws2.Range("A4").Value = "Intern verzoek / Actie voor medewerker"
$ws2.Range("B4").Value = 6
